$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reservasi")

# Insert two new columns before the existing "tanggal_check_in" / "tanggal_check_out"
# columns (G:H) to hold the new "jadwal_check_in" / "jadwal_check_out" (scheduled
# check-in/out) data. This shifts the old G:K -> I:M, preserving all existing
# values/formulas/styles in place.
$ws.Columns("G:H").Insert()

# New header labels for the inserted columns.
$ws.Range("G1").Value = "jadwal_check_in"
$ws.Range("H1").Value = "jadwal_check_out"

# Populate the new "jadwal_check_in" / "jadwal_check_out" values for every
# reservation row.
$ws.Range("G2").Value = "2022-04-25"
$ws.Range("H2").Value = "2022-04-27"

$ws.Range("G3").Value = "2022-04-28"
$ws.Range("H3").Value = "2022-05-01"

$ws.Range("G4").Value = "2022-04-20"
$ws.Range("H4").Value = "2022-04-24"

$ws.Range("G5").Value = "2022-04-22"
$ws.Range("H5").Value = "2022-04-26"

$ws.Range("G6").Value = "2022-04-24"
$ws.Range("H6").Value = "2022-04-29"

$ws.Range("G7").Value = "2022-04-27"
$ws.Range("H7").Value = "2022-04-29"

$ws.Range("G8").Value = "2022-04-26"
$ws.Range("H8").Value = "2022-05-18"

$ws.Range("G9").Value = "2022-05-06"
$ws.Range("H9").Value = "2022-05-09"

$ws.Range("G10").Value = "2022-04-30"
$ws.Range("H10").Value = "2022-05-05"

$ws.Range("G11").Value = "2022-05-11"
$ws.Range("H11").Value = "2022-05-14"

$ws.Range("G12").Value = "2022-05-13"
$ws.Range("H12").Value = "2022-05-18"

$ws.Range("G13").Value = "2022-05-17"
$ws.Range("H13").Value = "2022-05-18"

$ws.Range("G14").Value = "2022-05-16"
$ws.Range("H14").Value = "2022-05-21"

$ws.Range("G15").Value = "2022-05-16"
$ws.Range("H15").Value = "2022-05-21"

$ws.Range("G16").Value = "2022-05-16"
$ws.Range("H16").Value = "2022-05-21"

# Match the recalculated "best fit" column widths for the new columns as close
# as this engine's column-width model allows.
$ws.Columns("G:G").ColumnWidth = 14.251
$ws.Columns("H:H").ColumnWidth = 15.584
